$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.379369
$ws.Range("H2").Value = 31.138107
$ws.Range("I2").Value = 0.01614698522449884
$ws.Range("J2").Value = 0.01614698522449883
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 181.8232256666666
$ws.Range("N2").Value = 545.4696769999999
$ws.Range("O2").Value = 0.5898296910336229
$ws.Range("P2").Value = 0.5898296910336229
$ws.Range("Q2").Value = 1887.210351964604
$ws.Range("R2").Value = 16984.89316768144
$ws.Range("S2").Value = 0.009523971306090622
$ws.Range("T2").Value = 0.009523971306090621

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.379369
$ws.Range("H3").Value = 31.138107
$ws.Range("I3").Value = 0.01614698522449884
$ws.Range("J3").Value = 0.01614698522449883
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.367786666666667
$ws.Range("N3").Value = 7.10336
$ws.Range("O3").Value = 0.007681036748263821
$ws.Range("P3").Value = 0.007681036748263821
$ws.Range("Q3").Value = 24.57613152661334
$ws.Range("R3").Value = 221.18518373952
$ws.Range("S3").Value = 0.0001240255868830485
$ws.Range("T3").Value = 0.0001240255868830485

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.379369
$ws.Range("H4").Value = 31.138107
$ws.Range("I4").Value = 0.01614698522449884
$ws.Range("J4").Value = 0.01614698522449883
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.0729296666667
$ws.Range("N4").Value = 372.218789
$ws.Range("O4").Value = 0.4024892722181133
$ws.Range("P4").Value = 0.4024892722181133
$ws.Range("Q4").Value = 1287.79871992138
$ws.Range("R4").Value = 11590.18847929242
$ws.Range("S4").Value = 0.006498988331525165
$ws.Range("T4").Value = 0.006498988331525163

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 604.0312093333333
$ws.Range("H5").Value = 1812.093628
$ws.Range("I5").Value = 0.9396797639857967
$ws.Range("J5").Value = 0.9396797639857967
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 181.8232256666666
$ws.Range("N5").Value = 545.4696769999999
$ws.Range("O5").Value = 0.5898296910336229
$ws.Range("P5").Value = 0.5898296910336229
$ws.Range("Q5").Value = 109826.9028843242
$ws.Range("R5").Value = 988442.1259589181
$ws.Range("S5").Value = 0.5542510248622902
$ws.Range("T5").Value = 0.5542510248622902

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 604.0312093333333
$ws.Range("H6").Value = 1812.093628
$ws.Range("I6").Value = 0.9396797639857967
$ws.Range("J6").Value = 0.9396797639857967
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.367786666666667
$ws.Range("N6").Value = 7.10336
$ws.Range("O6").Value = 0.007681036748263821
$ws.Range("P6").Value = 0.007681036748263821
$ws.Range("Q6").Value = 1430.217043710009
$ws.Range("R6").Value = 12871.95339339008
$ws.Range("S6").Value = 0.007217714798774778
$ws.Range("T6").Value = 0.007217714798774778

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 604.0312093333333
$ws.Range("H7").Value = 1812.093628
$ws.Range("I7").Value = 0.9396797639857967
$ws.Range("J7").Value = 0.9396797639857967
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.0729296666667
$ws.Range("N7").Value = 372.218789
$ws.Range("O7").Value = 0.4024892722181133
$ws.Range("P7").Value = 0.4024892722181133
$ws.Range("Q7").Value = 74943.92175208627
$ws.Range("R7").Value = 674495.2957687766
$ws.Range("S7").Value = 0.3782110243247317
$ws.Range("T7").Value = 0.3782110243247317

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.39480333333333
$ws.Range("H8").Value = 85.18441
$ws.Range("I8").Value = 0.04417325078970442
$ws.Range("J8").Value = 0.04417325078970442
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 181.8232256666666
$ws.Range("N8").Value = 545.4696769999999
$ws.Range("O8").Value = 0.5898296910336229
$ws.Range("P8").Value = 0.5898296910336229
$ws.Range("Q8").Value = 5162.834734237284
$ws.Range("R8").Value = 46465.51260813556
$ws.Range("S8").Value = 0.0260546948652421
$ws.Range("T8").Value = 0.0260546948652421

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.39480333333333
$ws.Range("H9").Value = 85.18441
$ws.Range("I9").Value = 0.04417325078970442
$ws.Range("J9").Value = 0.04417325078970442
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.367786666666667
$ws.Range("N9").Value = 7.10336
$ws.Range("O9").Value = 0.007681036748263821
$ws.Range("P9").Value = 0.007681036748263821
$ws.Range("Q9").Value = 67.2328367352889
$ws.Range("R9").Value = 605.0955306176
$ws.Range("S9").Value = 0.0003392963626059935
$ws.Range("T9").Value = 0.0003392963626059935

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.39480333333333
$ws.Range("H10").Value = 85.18441
$ws.Range("I10").Value = 0.04417325078970442
$ws.Range("J10").Value = 0.04417325078970442
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.0729296666667
$ws.Range("N10").Value = 372.218789
$ws.Range("O10").Value = 0.4024892722181133
$ws.Range("P10").Value = 0.4024892722181133
$ws.Range("Q10").Value = 3523.026436875499
$ws.Range("R10").Value = 31707.23793187949
$ws.Range("S10").Value = 0.01777925956185633
$ws.Range("T10").Value = 0.01777925956185633
